# Update the "final main contribution" row labels: the commit renames the
# aggregate "Sharks" category labels so "Rays" / "Chimaeras" are lower-case
# ("rays" / "chimaeras") to match the rest of the sentence casing.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2:A3").Value2 = "Sharks (incl. rays & chimaeras)_national"
$ws.Range("A4:A5").Value2 = "Sharks (incl. rays & chimaeras)_global"

# Move the live selection/active cell (was K8, an empty cell far outside
# the used range) to A5, matching where the author was last looking when
# they saved this revision.
[void]$ws.Range("A5").Select()
